$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) counts
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 84
$wsExpo.Range("F8").Value = 125
$wsExpo.Range("F9").Value = 8951
$wsExpo.Range("F18").Value = 301
$wsExpo.Range("F19").Value = 73
$wsExpo.Range("F21").Value = 1117

# Sheet "全部类型" - same events, same field updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 84
$wsAll.Range("F10").Value = 125
$wsAll.Range("F11").Value = 8951
$wsAll.Range("F20").Value = 301
$wsAll.Range("F21").Value = 73
$wsAll.Range("F23").Value = 1117
